$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# NOTE on this runtime's worksheet/range handles: they are resolved by
# position, not stable object identity. Once a sheet is inserted/copied (or
# rows are inserted), any previously-captured $ws / $range variable will
# silently refer to "whatever is now at that position" rather than the
# original object. So below, every handle is (re)captured fresh, right after
# any operation that changes the sheet count or row/column layout, and is
# used immediately rather than being cached across such an operation.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new row for 2022-Q4 right after the
#    header row, pushing the existing quarters down by one row.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)
$total.Rows.Item(2).Insert()

# Copy the number-cell formatting (style) from the row below (which still
# carries the original "data row" style) so the new row matches the existing
# ones, then copy the "A column" label style (bold/bordered) separately.
$total.Range("B4:D4").Copy()
$total.Range("B2:D2").PasteSpecial(-4122)
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 8
$total.Range("D2").Value = 0.2

# ---------------------------------------------------------------------------
# 2) Add a brand-new "2022-Q4" worksheet right after "总计", built from a
#    copy of an existing quarterly sheet so it inherits the same sheet-level
#    properties (outline props, margins, header style, etc.), then its
#    contents are wiped and replaced with the 2022-Q4 holdings data.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item(2).Copy($null, $wb.Worksheets.Item(1))

# The new copy lands at index 2 ("总计" is 1); the sheet that was copied
# (the old "2022-Q3") is now pushed down to index 3 - use it as the style
# template for the new sheet before touching anything else.
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"
$tmpl = $wb.Worksheets.Item(3)

$q4.Cells.Clear()

# Header row style (bold / centered / bordered) copied from the template.
$tmpl.Range("B1:H1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Column A (row index) style, copied from the template's own data rows.
$tmpl.Range("A2").Copy()
$q4.Range("A2:A9").PasteSpecial(-4122)
$q4.Range("A2").Value = 0
$q4.Range("A3").Value = 1
$q4.Range("A4").Value = 2
$q4.Range("A5").Value = 3
$q4.Range("A6").Value = 4
$q4.Range("A7").Value = 5
$q4.Range("A8").Value = 6
$q4.Range("A9").Value = 7

# Fund code / fund name columns are plain text (codes like "000179" must
# keep their leading zeros), so force text format before assigning, then
# drop the formatting again so the cells end up unstyled, matching the rest
# of the workbook's data cells.
$q4.Range("B2:C9").NumberFormat = "@"
$q4.Range("B2").Value = "000179"
$q4.Range("C2").Value = "广发美国房地产指数（QDII）人民币A"
$q4.Range("B3").Value = "000180"
$q4.Range("C3").Value = "广发美国房地产指数（QDII）美元A"
$q4.Range("B4").Value = "160140"
$q4.Range("C4").Value = "南方道琼斯美国精选REIT指数（QDII-LOF）A"
$q4.Range("B5").Value = "070031"
$q4.Range("C5").Value = "嘉实全球房地产（QDII）"
$q4.Range("B6").Value = "160141"
$q4.Range("C6").Value = "南方道琼斯美国精选REIT指数（QDII-LOF）C"
$q4.Range("B7").Value = "320017"
$q4.Range("C7").Value = "诺安全球收益不动产（QDII）"
$q4.Range("B8").Value = "016278"
$q4.Range("C8").Value = "广发美国房地产指数（QDII）人民币C"
$q4.Range("B9").Value = "016279"
$q4.Range("C9").Value = "广发美国房地产指数（QDII）美元C"
$q4.Range("B2:C9").ClearFormats()

# Fund size / position / ratio / market-value columns are also stored as
# plain text in the source data (e.g. "1.60", "0.0606"), so the same
# text-then-clear-format trick is used here.
$q4.Range("D2:G9").NumberFormat = "@"
$q4.Range("D2").Value = "1.60"
$q4.Range("E2").Value = "92.49"
$q4.Range("F2").Value = "3.79"
$q4.Range("G2").Value = "0.0606"
$q4.Range("D3").Value = "1.60"
$q4.Range("E3").Value = "92.49"
$q4.Range("F3").Value = "3.79"
$q4.Range("G3").Value = "0.0606"
$q4.Range("D4").Value = "0.80"
$q4.Range("E4").Value = "92.31"
$q4.Range("F4").Value = "4.12"
$q4.Range("G4").Value = "0.0330"
$q4.Range("D5").Value = "0.39"
$q4.Range("E5").Value = "94.72"
$q4.Range("F5").Value = "4.65"
$q4.Range("G5").Value = "0.0181"
$q4.Range("D6").Value = "0.39"
$q4.Range("E6").Value = "92.31"
$q4.Range("F6").Value = "4.12"
$q4.Range("G6").Value = "0.0161"
$q4.Range("D7").Value = "0.24"
$q4.Range("E7").Value = "68.42"
$q4.Range("F7").Value = "5.60"
$q4.Range("G7").Value = "0.0134"
$q4.Range("D8").Value = "0.01"
$q4.Range("E8").Value = "92.49"
$q4.Range("F8").Value = "3.79"
$q4.Range("G8").Value = "0.0004"
$q4.Range("D9").Value = "0.01"
$q4.Range("E9").Value = "92.49"
$q4.Range("F9").Value = "3.79"
$q4.Range("G9").Value = "0.0004"
$q4.Range("D2:G9").ClearFormats()

# Position-rank column is a real number.
$q4.Range("H2").Value = 4
$q4.Range("H3").Value = 4
$q4.Range("H4").Value = 4
$q4.Range("H5").Value = 2
$q4.Range("H6").Value = 4
$q4.Range("H7").Value = 4
$q4.Range("H8").Value = 4
$q4.Range("H9").Value = 4
